$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the address at A3 (was "BUKIT BATOK FLYOVER") to the new value
$ws.Range("A3").Value = "BUKIT BATOK FIRE STATION"

# Update selection to match the new active cell / selection range
$ws.Range("A3").Select()
